# Rename the price-list worksheet, appending "(V)" to mark it as the
# viewer/web version, e.g. "Artikelpreisliste_11-2025" -> "Artikelpreisliste_11-2025(V)".
$wb = $excel.ActiveWorkbook

$oldName = "Artikelpreisliste_11-2025"
$newName = "Artikelpreisliste_11-2025(V)"

$ws = $wb.Worksheets.Item($oldName)
$ws.Name = $newName

# The sheet-scoped built-in "Print_Area" name embeds the (now stale) sheet
# name in its formula; re-point it explicitly so it follows the rename the
# same way the "_FilterDatabase" name already does.
$ws.PageSetup.PrintArea = '$A$1:$G$850'
